$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two report-section headings (B7 first so shared-string
# table order matches: "2. ..." gets interned before "1. ...")
$ws.Range("B7").Value = "2. Lượt truy cập các trang chi tiết sản phẩm:"
$ws.Range("B6").Value = "1. Số lượt truy cập trang chủ:"

# Move the active selection cursor to K19
$ws.Range("K19").Select()
